# Regenerate the "K" (strike) column (G) values on Sheet1, rows 2-10.
# This mirrors a re-run of the save_data generation process that now derives
# strike counts (K) directly rather than from "Strike#", and also
# recalculates std/mean and writes the resulting s_vals into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K values computed by the regenerated save_data routine, keyed by row.
$newK = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
